# Re-run of the DD-trigger optimization sweep: rows 2-27 get the fresh
# results (the previous run's stale tail, old rows 28-31, is dropped), and
# the one-off highlighted "best result" row loses its special formatting
# since the new winner no longer needs to be called out that way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns: r, E=DD_thres, G=Final_equity, H=Accs_started, I=Accs_alive, J=Accs_blown, K=Portf_max_DD
$data = @(
    @(2, 500, 187932, 15, 15, 0, -14056),
    @(3, 600, 143067, 11, 11, 0, -9613),
    @(4, 700, 101427, 8, 8, 0, -6268),
    @(5, 800, 81792, 7, 7, 0, -5195),
    @(6, 1000, 53612, 4, 4, 0, -3220),
    @(7, 900, 53295, 4, 4, 0, -2994),
    @(8, 1100, 32228, 2, 2, 0, -2147),
    @(9, 2100, 22240, 1, 1, 0, -1136),
    @(10, 2900, 22240, 1, 1, 0, -1136),
    @(11, 2800, 22240, 1, 1, 0, -1136),
    @(12, 2700, 22240, 1, 1, 0, -1136),
    @(13, 2600, 22240, 1, 1, 0, -1136),
    @(14, 2500, 22240, 1, 1, 0, -1136),
    @(15, 2400, 22240, 1, 1, 0, -1136),
    @(16, 2300, 22240, 1, 1, 0, -1136),
    @(17, 2200, 22240, 1, 1, 0, -1136),
    @(18, 1800, 22240, 1, 1, 0, -1136),
    @(19, 2000, 22240, 1, 1, 0, -1136),
    @(20, 1900, 22240, 1, 1, 0, -1136),
    @(21, 1700, 22240, 1, 1, 0, -1136),
    @(22, 1600, 22240, 1, 1, 0, -1136),
    @(23, 1500, 22240, 1, 1, 0, -1136),
    @(24, 1400, 22240, 1, 1, 0, -1136),
    @(25, 1300, 22240, 1, 1, 0, -1136),
    @(26, 1200, 22240, 1, 1, 0, -1136),
    @(27, 3000, 22240, 1, 1, 0, -1136)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $ws.Cells.Item($r, 9).Value = $row[4]
    $ws.Cells.Item($r, 10).Value = $row[5]
    $ws.Cells.Item($r, 11).Value = $row[6]
}

# Row 23 used to be specially bolded/highlighted (yellow fill) to flag a
# standout result; after the re-run no row needs that treatment, so strip
# it back to the plain style used everywhere else (keeping just the euro
# number format on the Final_equity cell, like every other data row).
$ws.Rows(23).ClearFormats()
$ws.Range("G23").NumberFormat = '#,##0.00\ "€"'

# The old sweep had four extra (now obsolete) result rows at the bottom.
$ws.Rows("28:31").Delete()

# Move the view: no more scrolled/frozen top-left cell, and the selection
# moves to B6.
$ws.Range("B6").Select()

# Match the saved window size from the edited workbook.
$win = $excel.ActiveWindow
$win.Width = 16200
$win.Height = 24825
